# Fruta / hortaliza, semanal
# Insert two new weekly price records (rows 212-213) into the Palta
# (avocado) price sheet for "Vega Monumental Concepción". Existing rows
# 212-305 shift down to 214-307; the sheet's used range grows from
# A1:T305 to A1:T307.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 212 (shifts everything
# below down by two rows and carries the row-above's formatting, which
# is what gives the new D212/D213 cells the existing date number format).
$ws.Rows("212:213").Insert()

# --- Row 212: Hass / Primera, Perú, $/bandeja 10 kilos ---
$ws.Cells.Item(212, 1).Value = 11
$ws.Cells.Item(212, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(212, 3).Value = "Bíobío"
$ws.Cells.Item(212, 4).Value = 44468
$ws.Cells.Item(212, 5).Value = 8
$ws.Cells.Item(212, 6).Value = "Fruta"
$ws.Cells.Item(212, 7).Value = 100106
$ws.Cells.Item(212, 8).Value = "Oleaginosos"
$ws.Cells.Item(212, 9).Value = 100106002
$ws.Cells.Item(212, 10).Value = "Palta"
$ws.Cells.Item(212, 11).Value = "Hass"
$ws.Cells.Item(212, 12).Value = "Primera"
$ws.Cells.Item(212, 13).Value = 100
$ws.Cells.Item(212, 14).Value = 26000
$ws.Cells.Item(212, 15).Value = 27000
$ws.Cells.Item(212, 16).Value = 26500
$ws.Cells.Item(212, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(212, 18).Value = "Perú"
$ws.Cells.Item(212, 19).Value = 2650
$ws.Cells.Item(212, 20).Value = 10

# --- Row 213: Hass / Segunda, Perú, $/bandeja 10 kilos ---
$ws.Cells.Item(213, 1).Value = 11
$ws.Cells.Item(213, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(213, 3).Value = "Bíobío"
$ws.Cells.Item(213, 4).Value = 44468
$ws.Cells.Item(213, 5).Value = 8
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100106
$ws.Cells.Item(213, 8).Value = "Oleaginosos"
$ws.Cells.Item(213, 9).Value = 100106002
$ws.Cells.Item(213, 10).Value = "Palta"
$ws.Cells.Item(213, 11).Value = "Hass"
$ws.Cells.Item(213, 12).Value = "Segunda"
$ws.Cells.Item(213, 13).Value = 50
$ws.Cells.Item(213, 14).Value = 24000
$ws.Cells.Item(213, 15).Value = 24000
$ws.Cells.Item(213, 16).Value = 24000
$ws.Cells.Item(213, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(213, 18).Value = "Perú"
$ws.Cells.Item(213, 19).Value = 2400
$ws.Cells.Item(213, 20).Value = 10
